# Add a reviewer comment ("Vic4ever") anchored to the "Server Information
# will be collect from agents." sentence in the Data Flows section, per the
# commit "Add questions and comments to review".

$d = $word.ActiveDocument

# Locate the sentence the comment is anchored to.
$commentRange = $d.Content
$found = $commentRange.Find.Execute(
    "Server Information will be collect from agents.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $comment = $d.Comments.Add($commentRange, "Có cần thêm 1 cái để quét server info ?")
    $comment.Author = "Vic4ever"
    $comment.Initial = "V"
}
